# Update "想去人数" (F column) counts across sheets, matching the
# gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13
$ws1.Range("F4").Value = 93
$ws1.Range("F5").Value = 911
$ws1.Range("F7").Value = 7134
$ws1.Range("F9").Value = 156
$ws1.Range("F10").Value = 6551
$ws1.Range("F13").Value = 4519
$ws1.Range("F17").Value = 4579
$ws1.Range("F19").Value = 249
$ws1.Range("F20").Value = 20
$ws1.Range("F21").Value = 347
$ws1.Range("F28").Value = 8194
$ws1.Range("F30").Value = 1420
$ws1.Range("F32").Value = 718
$ws1.Range("F34").Value = 50
$ws1.Range("F37").Value = 1679
$ws1.Range("F40").Value = 43
$ws1.Range("F41").Value = 4241
$ws1.Range("F44").Value = 118
$ws1.Range("F47").Value = 1126
$ws1.Range("F48").Value = 5
$ws1.Range("F49").Value = 22

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 32

# --- Sheet "本地生活" --- (no changes)

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 13
$ws4.Range("F7").Value = 93
$ws4.Range("F9").Value = 911
$ws4.Range("F11").Value = 7134
$ws4.Range("F13").Value = 156
$ws4.Range("F14").Value = 6551
$ws4.Range("F17").Value = 4520
$ws4.Range("F21").Value = 4580
$ws4.Range("F23").Value = 249
$ws4.Range("F24").Value = 347
$ws4.Range("F29").Value = 8194
$ws4.Range("F31").Value = 1420
$ws4.Range("F33").Value = 718
$ws4.Range("F35").Value = 50
$ws4.Range("F37").Value = 1679
$ws4.Range("F41").Value = 4241
$ws4.Range("F44").Value = 118
$ws4.Range("F47").Value = 1126
$ws4.Range("F48").Value = 5
$ws4.Range("F49").Value = 22
